# Fruta / hortaliza, semanal
# Insert a new weekly record for "Vega Monumental Concepción" (Naranja)
# right after the existing row 140, shifting the rest of the price-history
# rows (old 141-170) down by one (new 142-171).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 141, pushing old rows 141:170 down to 142:171
$ws.Rows("141:141").Insert()

# Populate the new row with the latest weekly observation
$ws.Cells.Item(141, 1).Value = 11
$ws.Cells.Item(141, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(141, 3).Value = "Bíobío"
$ws.Cells.Item(141, 4).Value = 44511
$ws.Cells.Item(141, 5).Value = 8
$ws.Cells.Item(141, 6).Value = "Fruta"
$ws.Cells.Item(141, 7).Value = 100102
$ws.Cells.Item(141, 8).Value = "Cítricos"
$ws.Cells.Item(141, 9).Value = 100102005
$ws.Cells.Item(141, 10).Value = "Naranja"
$ws.Cells.Item(141, 11).Value = "Navel Late"
$ws.Cells.Item(141, 12).Value = "Primera"
$ws.Cells.Item(141, 13).Value = 200
$ws.Cells.Item(141, 14).Value = 7500
$ws.Cells.Item(141, 15).Value = 8000
$ws.Cells.Item(141, 16).Value = 7750
$ws.Cells.Item(141, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(141, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(141, 19).Value = 517
$ws.Cells.Item(141, 20).Value = 15

# Keep the date column's number format consistent with the rest of the sheet
$ws.Cells.Item(141, 4).NumberFormat = $ws.Cells.Item(142, 4).NumberFormat
